$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Logs")
$ws2 = $wb.Worksheets.Item("Weekly Summary")

# --- Fill in Logs rows 57-75 with Week 5 worklog entries ---
$ws1.Range("A57").Value2 = 'Tithra Chap'
$ws1.Range("B57").Value2 = 'As. Prof. Richard Dazeley'
$ws1.Range("C57").Value2 = 'Emotion Recognition Using Facial Expression'
$ws1.Range("D57").Value2 = 44417
$ws1.Range("E57").Value2 = 60
$ws1.Range("F57").Value2 = 'Supervision Meeting'
$ws1.Range("G57").Value2 = 'Week 5 meeting with superviors'

$ws1.Range("A58").Value2 = 'Tithra Chap'
$ws1.Range("B58").Value2 = 'As. Prof. Richard Dazeley'
$ws1.Range("C58").Value2 = 'Emotion Recognition Using Facial Expression'
$ws1.Range("D58").Value2 = 44417
$ws1.Range("E58").Value2 = 120
$ws1.Range("F58").Value2 = 'Artefact Implementation'
$ws1.Range("G58").Value2 = 'Integrate the MobileNet with FER2013'

$ws1.Range("A59").Value2 = 'Tithra Chap'
$ws1.Range("B59").Value2 = 'As. Prof. Richard Dazeley'
$ws1.Range("C59").Value2 = 'Emotion Recognition Using Facial Expression'
$ws1.Range("D59").Value2 = 44417
$ws1.Range("E59").Value2 = 60
$ws1.Range("F59").Value2 = 'Artefact Implementation'
$ws1.Range("G59").Value2 = 'Integrate the MobileNet with FER2013'

$ws1.Range("A60").Value2 = 'Tithra Chap'
$ws1.Range("B60").Value2 = 'As. Prof. Richard Dazeley'
$ws1.Range("C60").Value2 = 'Emotion Recognition Using Facial Expression'
$ws1.Range("D60").Value2 = 44418
$ws1.Range("E60").Value2 = 120
$ws1.Range("F60").Value2 = 'Artefact Implementation'
$ws1.Range("G60").Value2 = 'Hypertune the experiment of MobileNet with FER2013 and generate results'

$ws1.Range("A61").Value2 = 'Tithra Chap'
$ws1.Range("B61").Value2 = 'As. Prof. Richard Dazeley'
$ws1.Range("C61").Value2 = 'Emotion Recognition Using Facial Expression'
$ws1.Range("D61").Value2 = 44418
$ws1.Range("E61").Value2 = 120
$ws1.Range("F61").Value2 = 'Artefact Implementation'
$ws1.Range("G61").Value2 = 'Hypertune the experiment of MobileNet with FER2013 and generate results'

$ws1.Range("A62").Value2 = 'Tithra Chap'
$ws1.Range("B62").Value2 = 'As. Prof. Richard Dazeley'
$ws1.Range("C62").Value2 = 'Emotion Recognition Using Facial Expression'
$ws1.Range("D62").Value2 = 44418
$ws1.Range("E62").Value2 = 80
$ws1.Range("F62").Value2 = 'Artefact Implementation'
$ws1.Range("G62").Value2 = 'Hypertune the experiment of MobileNet with FER2013 and generate results'

$ws1.Range("A63").Value2 = 'Tithra Chap'
$ws1.Range("B63").Value2 = 'As. Prof. Richard Dazeley'
$ws1.Range("C63").Value2 = 'Emotion Recognition Using Facial Expression'
$ws1.Range("D63").Value2 = 44419
$ws1.Range("E63").Value2 = 120
$ws1.Range("F63").Value2 = 'Artefact Implementation'
$ws1.Range("G63").Value2 = 'Hypertune the experiment of MobileNet with FER2013 and generate results'

$ws1.Range("A64").Value2 = 'Tithra Chap'
$ws1.Range("B64").Value2 = 'As. Prof. Richard Dazeley'
$ws1.Range("C64").Value2 = 'Emotion Recognition Using Facial Expression'
$ws1.Range("D64").Value2 = 44419
$ws1.Range("E64").Value2 = 120
$ws1.Range("F64").Value2 = 'Artefact Implementation'
$ws1.Range("G64").Value2 = 'Hypertune the experiment of MobileNet with FER2013 and generate results'

$ws1.Range("A65").Value2 = 'Tithra Chap'
$ws1.Range("B65").Value2 = 'As. Prof. Richard Dazeley'
$ws1.Range("C65").Value2 = 'Emotion Recognition Using Facial Expression'
$ws1.Range("D65").Value2 = 44419
$ws1.Range("E65").Value2 = 120
$ws1.Range("F65").Value2 = 'Artefact Implementation'
$ws1.Range("G65").Value2 = 'Hypertune the experiment of MobileNet with FER2013 and generate results'

$ws1.Range("A66").Value2 = 'Tithra Chap'
$ws1.Range("B66").Value2 = 'As. Prof. Richard Dazeley'
$ws1.Range("C66").Value2 = 'Emotion Recognition Using Facial Expression'
$ws1.Range("D66").Value2 = 44421
$ws1.Range("E66").Value2 = 120
$ws1.Range("F66").Value2 = 'Literature Search'
$ws1.Range("G66").Value2 = 'Find pre-processing and parameter tuning techniques in good papers'

$ws1.Range("A67").Value2 = 'Tithra Chap'
$ws1.Range("B67").Value2 = 'As. Prof. Richard Dazeley'
$ws1.Range("C67").Value2 = 'Emotion Recognition Using Facial Expression'
$ws1.Range("D67").Value2 = 44421
$ws1.Range("E67").Value2 = 120
$ws1.Range("F67").Value2 = 'Literature Search'
$ws1.Range("G67").Value2 = 'Find pre-processing and parameter tuning techniques in good papers'

$ws1.Range("A68").Value2 = 'Tithra Chap'
$ws1.Range("B68").Value2 = 'As. Prof. Richard Dazeley'
$ws1.Range("C68").Value2 = 'Emotion Recognition Using Facial Expression'
$ws1.Range("D68").Value2 = 44421
$ws1.Range("E68").Value2 = 90
$ws1.Range("F68").Value2 = 'Literature Search'
$ws1.Range("G68").Value2 = 'Find pre-processing and parameter tuning techniques in good papers'

$ws1.Range("A69").Value2 = 'Tithra Chap'
$ws1.Range("B69").Value2 = 'As. Prof. Richard Dazeley'
$ws1.Range("C69").Value2 = 'Emotion Recognition Using Facial Expression'
$ws1.Range("D69").Value2 = 44422
$ws1.Range("E69").Value2 = 120
$ws1.Range("F69").Value2 = 'Artefact Implementation'
$ws1.Range("G69").Value2 = 'Locate CK+ dataset and manipluate the preprocessing'

$ws1.Range("A70").Value2 = 'Tithra Chap'
$ws1.Range("B70").Value2 = 'As. Prof. Richard Dazeley'
$ws1.Range("C70").Value2 = 'Emotion Recognition Using Facial Expression'
$ws1.Range("D70").Value2 = 44422
$ws1.Range("E70").Value2 = 120
$ws1.Range("F70").Value2 = 'Artefact Implementation'
$ws1.Range("G70").Value2 = 'Locate CK+ dataset and manipluate the preprocessing'

$ws1.Range("A71").Value2 = 'Tithra Chap'
$ws1.Range("B71").Value2 = 'As. Prof. Richard Dazeley'
$ws1.Range("C71").Value2 = 'Emotion Recognition Using Facial Expression'
$ws1.Range("D71").Value2 = 44422
$ws1.Range("E71").Value2 = 40
$ws1.Range("F71").Value2 = 'Artefact Implementation'
$ws1.Range("G71").Value2 = 'Locate CK+ dataset and manipluate the preprocessing'

$ws1.Range("A72").Value2 = 'Tithra Chap'
$ws1.Range("B72").Value2 = 'As. Prof. Richard Dazeley'
$ws1.Range("C72").Value2 = 'Emotion Recognition Using Facial Expression'
$ws1.Range("D72").Value2 = 44423
$ws1.Range("E72").Value2 = 120
$ws1.Range("F72").Value2 = 'Artefact Implementation'
$ws1.Range("G72").Value2 = 'Experiment ResNet20 with CK+ dataset'

$ws1.Range("A73").Value2 = 'Tithra Chap'
$ws1.Range("B73").Value2 = 'As. Prof. Richard Dazeley'
$ws1.Range("C73").Value2 = 'Emotion Recognition Using Facial Expression'
$ws1.Range("D73").Value2 = 44423
$ws1.Range("E73").Value2 = 120
$ws1.Range("F73").Value2 = 'Artefact Implementation'
$ws1.Range("G73").Value2 = 'Experiment ResNet20 with CK+ dataset'

$ws1.Range("A74").Value2 = 'Tithra Chap'
$ws1.Range("B74").Value2 = 'As. Prof. Richard Dazeley'
$ws1.Range("C74").Value2 = 'Emotion Recognition Using Facial Expression'
$ws1.Range("D74").Value2 = 44423
$ws1.Range("E74").Value2 = 110
$ws1.Range("F74").Value2 = 'Artefact Implementation'
$ws1.Range("G74").Value2 = 'Experiment ResNet20 with CK+ dataset'

$ws1.Range("A75").Value2 = 'Tithra Chap'
$ws1.Range("B75").Value2 = 'As. Prof. Richard Dazeley'
$ws1.Range("C75").Value2 = 'Emotion Recognition Using Facial Expression'
$ws1.Range("D75").Value2 = 44423
$ws1.Range("E75").Value2 = 120
$ws1.Range("F75").Value2 = 'OnTrack Task'
$ws1.Range("G75").Value2 = 'Work task 5.1P'

# --- Extend the Logs sheet with a new formatted blank row 250 (copy date format from D249) ---
$ws1.Range("D249").Copy() | Out-Null
$ws1.Range("D250").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Weekly Summary: Week 5 totals (row 7) ---
$ws2.Range("B7").Formula = "=ROUNDDOWN(SUM(Logs!E57:E75)/60,0)"
$ws2.Range("C7").Formula = "=MOD(SUM(Logs!E57:E75),60)"

# --- Weekly Summary: Total row (row 15) now spans one additional row ---
$ws2.Range("B15").Formula = "=ROUNDDOWN(SUM(Logs!E3:E201)/60,0)"
$ws2.Range("C15").Formula = "=MOD(SUM(Logs!E3:E201),60)"

# --- View state: Weekly Summary selection, then make Logs the active/selected sheet ---
$ws2.Activate()
$ws2.Range("C16").Select() | Out-Null

$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 64
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("C76").Select() | Out-Null
